$d = $word.ActiveDocument

# 1. Remove the stray "Interface gráfica..." caption prefix from the
#    intimation paragraph, keeping the rest of the sentence intact.
$d.Content.Find.Execute(
    "Interface gráfica do usuário, Aplicativo, Word Descrição gerada automaticamente. Tendo em vista",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tendo em vista", 2)

# 2. Company name
$d.Content.Find.Execute(
    "Equatorial Maranhão Distribuidora de Energia S.A.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Vivo S.A", 2)

# 3. Activity / segment
$d.Content.Find.Execute(
    "Distribuição de Energia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Operadora de telefonia", 2)

# 4. CNPJ
$d.Content.Find.Execute(
    "06.272.793/0001-84",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "06.952.209/344", 2)

# 5. Address (street)
$d.Content.Find.Execute(
    "AL A, Quadra SQS",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Avenida Engenheiro Luiz Carlos Berrini", 2)

# 6. Bairro
$d.Content.Find.Execute(
    "Loteamento Quintandinha Altos do Calhau",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cidade Monções", 2)

# 7. Cidade
$d.Content.Find.Execute(
    "São Luiz",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "São Paulo", 2)

# 9. CEP
$d.Content.Find.Execute(
    "65.070-900",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "04571-936", 2)

# 10. Fiscal agent name
$d.Content.Find.Execute(
    "Jacia Andréia Nascimento Sousa Pedral",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "André Luis da Silva Oliveira", 2)

# 11. Matrícula
$d.Content.Find.Execute(
    "996552-1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "734424", 2)

# 8/12/13. The ESTADO value cell only contains "MA" (too short/common a
# string to safely Find/Replace document-wide -- it also occurs inside
# unrelated sentences such as "São José de Ribamar MA"), and the "Nº" /
# "COMPLEMENTO" value cells are empty paragraphs that need a new run
# added. Tables(...).Cell(...).Range.Paragraphs collections are not
# reliably indexable in this runtime, so address the three paragraphs
# directly through the document-level Paragraphs collection instead,
# which keeps the edits scoped to exactly the right paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $style = $p.Style.NameLocal
    # Range.Text includes the trailing paragraph mark (and, for a
    # paragraph that ends a table cell, the cell-end mark too) so strip
    # those control characters before comparing. (TrimEnd needs explicit
    # [char] arguments here -- passing them as strings silently no-ops.)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($style -eq "12-bold" -and $text -eq "Nº") {
        $d.Paragraphs.Item($i + 1).Range.InsertAfter("928")
    }
    elseif ($style -eq "12-bold" -and $text -eq "COMPLEMENTO") {
        $d.Paragraphs.Item($i + 1).Range.InsertAfter("Parque Alvorada")
    }
    elseif ($style -eq "principal" -and $text -eq "MA") {
        $p.Range.Text = "SP"
    }
}
